# Append Trade #14 (closed) as a new row to the "leadlag" worksheet of the
# live trading results workbook. Matches the diff: a new row 13 is added
# with columns A:N populated, bumping the sheet's used-range dimension
# from A1:N12 to A1:N13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 13

$ws.Cells.Item($row, 1).Value = 14                 # A: Trade #
# B: Date - force text formatting first so the "YYYY-MM-DD" string is
# stored as literal text (matching every other row) instead of being
# auto-parsed into a date serial number; then restore the default
# "Normal" style so no stray per-cell formatting is left behind.
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "2026-02-16"
$ws.Cells.Item($row, 2).Style = "Normal"
$ws.Cells.Item($row, 3).Value = "21:23:08"         # C: Time
$ws.Cells.Item($row, 4).Value = "leadlag"          # D: Strategy
$ws.Cells.Item($row, 5).Value = "DOWN"             # E: Side
$ws.Cells.Item($row, 6).Value = 69382.565          # F: Entry Price
# G: Exit Price - still open, left blank
$ws.Cells.Item($row, 8).Value = "OPEN"             # H: Status
$ws.Cells.Item($row, 9).Value = 0                  # I: P&L %
$ws.Cells.Item($row, 10).Value = 0                 # J: P&L $
$ws.Cells.Item($row, 11).Value = 0.6303            # K: Confidence
$ws.Cells.Item($row, 12).Value = "Binance leading with -0.063% move"  # L: Entry Reason
# M: Exit Reason - still open, left blank
$ws.Cells.Item($row, 14).Value = 0                 # N: Duration (min)
